# Purchase report update: SALES REPORT -> PURCHASE STATUS REPORT
# Adds a "Status" column, updates the generated-on date, adds two new
# purchase rows and updates the total amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper for values that look like numbers/dates (e.g. "10000.00",
# "2019-01-18") so Excel keeps them as plain text instead of silently
# converting them to a number or date serial.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Title & generated-on date ---------------------------------------------
$ws.Cells.Item(1, 1).Value = "PURCHASE STATUS REPORT"
$ws.Cells.Item(2, 2).Value = "17/1/2019"

# --- New "Status" column header ---------------------------------------------
$ws.Cells.Item(4, 6).Value = "Status"

# --- Status values for the existing purchase rows (6-13) --------------------
$ws.Cells.Item(6, 6).Value = "purchased_done"
$ws.Cells.Item(9, 6).Value = "purchased_done"

# --- Insert two new rows before the "Total Amount" row (currently row 14) ---
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# --- Row 14: 9th purchase order ----------------------------------------------
$ws.Cells.Item(14, 1).Value = 9
$ws.Cells.Item(14, 2).Value = "po/abc/2019/5"
$ws.Cells.Item(14, 3).Value = "abc"
Set-TextCell 14 4 "2019-01-18"
Set-TextCell 14 5 "15000.00"

# --- Row 15: 10th purchase order ---------------------------------------------
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "po/abc/2019/6"
$ws.Cells.Item(15, 3).Value = "abc"
Set-TextCell 15 4 "2019-01-12"
Set-TextCell 15 5 "7500.00"

# --- Updated total (now on row 16 after the insert) ---------------------------
Set-TextCell 16 5 "285473.11"

# --- Match the author's final selection ---------------------------------------
$ws.Range("F4").Select()
